$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.978.81'
$ws.Range("E2").Value = '  +1.94%  '
$ws.Range("D3").Value = '2.528.59'
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("D9").Value = '2.530.30'
$ws.Range("E9").Value = '  +1.61%  '
$ws.Range("E10").Value = '  +7.75%  '
$ws.Range("E11").Value = '  -0.99%  '
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").Value = '2.990.50'
$ws.Range("E14").Value = '  +2.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.92%  '
$ws.Range("D16").Value = '68.890.33'
$ws.Range("E16").Value = '  +2.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000174'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.54%  '
$ws.Range("D18").Value = '2.541.44'
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '361.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("B26").Value = 'SuiNetwork'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.68'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.52%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.43%  '
$ws.Range("D28").Value = '2.657.47'
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '511.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").Value = '0.0₃0885'
$ws.Range("E31").Value = '  -2.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.78'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.77'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.89%  '
$ws.Range("E37").Value = '  -3.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.69%  '
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.31'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.79%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.59%  '
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.80'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("E45").Value = '  -2.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '151.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.38%  '
$ws.Range("E47").Value = '  +2.03%  '
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0742'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0251'
$ws.Range("E50").Value = '  -2.51%  '
$ws.Range("B51").Value = 'Optimism'
$ws.Range("C51").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.00%  '
